$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial for every data row (rows 2-90).
# The value was bumped by one day (45202 -> 45203, i.e. 2023-10-03 -> 2023-10-04).
$ws.Range("C2:C90").Value = 45203
